# Better handle copying paragraph styles:
# Whenever a "MSC_Join" (pStyle MSCJoin) paragraph is immediately
# followed by two plain "Normal" paragraphs (the "[...]" placeholder
# paragraph and the blank paragraph after it), propagate the MSCJoin
# style onto those two paragraphs as well, so the whole join/placeholder
# group shares consistent paragraph formatting.

$d = $word.ActiveDocument

$targets = New-Object System.Collections.ArrayList

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count - 2; $i++) {
    $anchor = $d.Paragraphs($i)
    if ($anchor.Style.NameLocal -eq "MSC_Join") {
        $next1 = $d.Paragraphs($i + 1)
        $next2 = $d.Paragraphs($i + 2)
        if (($next1.Style.NameLocal -eq "Normal") -and ($next2.Style.NameLocal -eq "Normal")) {
            [void]$targets.Add($i + 1)
            [void]$targets.Add($i + 2)
        }
    }
}

foreach ($idx in $targets) {
    $d.Paragraphs($idx).Style = "MSC_Join"
}
